# Adds answers for problem 2 and the first half of problem 3 to the
# Process Survey document, right after the existing problem 1 answer
# and before the section properties mark.

$d = $word.ActiveDocument

# Collapsed range sitting just before the document's final paragraph
# mark (i.e. right after the last existing paragraph's content).
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="style0"/></w:pPr><w:r><w:rPr/><w:t>2. a.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="style0"/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">b. ps aux | grep ^haley | awk '{ print $11 }' &gt; </w:t><w:tab/><w:t>~/Documents/gitRepos/cmsi387/homework/shell-gymnastics/answer2.txt</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="style0"/></w:pPr><w:r><w:rPr/><w:t>3. a. Real Memory: /usr/sbin/console-kit-daemon</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="style0"/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>b. used &#8220;man ps&#8221; as well</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="style0"/></w:pPr><w:r><w:rPr/><w:tab/><w:t xml:space="preserve">ps aux | awk '{print $5"\t"$11}' | sort -n -r &gt; </w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:t>~/Documents/gitRepos/cmsi387/homework/shell-gymnastics/answer3real.txt</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newParagraphsXml)
